# feat: add 2022-Q3 data
#
# Insert a new "2022-Q3" worksheet (fund holdings for 688567) right after the
# "总计" summary sheet, i.e. immediately before the existing "2022-Q2" sheet,
# and record the new quarter's totals at the top of the "总计" table.

$wb = $excel.ActiveWorkbook

# --- 1. Duplicate the "2022-Q2" sheet (it already carries the right column
#        layout/styles) and place the copy right before it; this becomes our
#        "2022-Q3" sheet. ---
$q2Sheet = $wb.Worksheets.Item(2)
$q2Sheet.Copy($q2Sheet)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# The template had 7 data rows (rows 2-8); the new sheet only needs 4
# (rows 2-5), so drop the trailing rows.
$newSheet.Rows.Item(8).Delete()
$newSheet.Rows.Item(7).Delete()
$newSheet.Rows.Item(6).Delete()

# --- 2. Overwrite the fund-holding data rows. Row 1 (headers) and column A
#        (row index 0..3) already match, so only B:H need new values.
#        Fund code (B) and columns D:G hold numeric-looking text
#        (e.g. "020026", "6.11", "0.70") so they must be written with a text
#        number-format to keep their literal formatting (leading zeros,
#        trailing zeros). Column C (fund name) is plain Chinese text and
#        stays text on its own, so it's left at the default format. ---
$textFormatCols = 2, 4, 5, 6, 7
$rows = @(
    @(2, "020026", "国泰成长优选混合", "6.11", "93.49", "3.31", "0.2022", 10),
    @(3, "014575", "鑫元清洁能源混合C", "1.77", "93.05", "5.81", "0.1028", 9),
    @(4, "014574", "鑫元清洁能源混合A", "0.82", "93.05", "5.81", "0.0476", 9),
    @(5, "233001", "大摩基础行业混合", "0.70", "78.50", "5.80", "0.0406", 7)
)

foreach ($row in $rows) {
    $r = $row[0]
    foreach ($c in $textFormatCols) {
        $newSheet.Cells.Item($r, $c).NumberFormat = "@"
        $newSheet.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
}

# --- 3. Insert the corresponding summary row at the top of "总计"'s table,
#        pushing the existing quarters down by one row. Shift content
#        manually (bottom-up) instead of Rows.Insert() so no stray
#        formatting gets pulled onto the shifted cells. ---
$summary = $wb.Worksheets.Item(1)

for ($r = 7; $r -ge 2; $r--) {
    $summary.Cells.Item($r + 1, 2).Value = $summary.Cells.Item($r, 2).Value2
    $summary.Cells.Item($r + 1, 3).Value = $summary.Cells.Item($r, 3).Value2
    $summary.Cells.Item($r + 1, 4).Value = $summary.Cells.Item($r, 4).Value2
}

# Column A is just the sequential row index (0-based) - keep it that way,
# and make sure the brand new row 8 picks up the same "index" style as the
# rest of column A.
$summary.Cells.Item(7, 1).Copy()
$summary.Cells.Item(8, 1).PasteSpecial(-4122)
for ($r = 2; $r -le 8; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}

$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 4
$summary.Cells.Item(2, 4).Value = 0.39
